$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the weekday/time header labels in row 1 (D1:W1): "Mon-8" -> "Mon08", etc.
$ws.Range("D1").Value = "Mon08"
$ws.Range("E1").Value = "Mon010"
$ws.Range("F1").Value = "Mon012"
$ws.Range("G1").Value = "Mon02"
$ws.Range("H1").Value = "Tue08"
$ws.Range("I1").Value = "Tue010"
$ws.Range("J1").Value = "Tue012"
$ws.Range("K1").Value = "Tue02"
$ws.Range("L1").Value = "Wed08"
$ws.Range("M1").Value = "Wed010"
$ws.Range("N1").Value = "Wed012"
$ws.Range("O1").Value = "Wed02"
$ws.Range("P1").Value = "Thu08"
$ws.Range("Q1").Value = "Thu010"
$ws.Range("R1").Value = "Thu012"
$ws.Range("S1").Value = "Thu02"
$ws.Range("T1").Value = "Fri08"
$ws.Range("U1").Value = "Fri010"
$ws.Range("V1").Value = "Fri012"
$ws.Range("W1").Value = "Fri02"

# D11 held the text placeholder "-" for missing data; replace with numeric 0.
$ws.Range("D11").Value = 0

# Move the active selection to I13.
$ws.Range("I13").Select()
